# Edit script for nservicebus/azure-service-bus/Scopes.pptx
# Summary of changes:
#  1. Duplicate slide 2 (the "Pipeline" diagram) to create a new slide,
#     inserted right after it (becomes slide 3), then rework its shapes
#     into the "Incoming/Outgoing/Handler Pipeline" scope diagram.
#  2. Rename the three "Pipeline" labels on (the original) slide 2 and
#     nudge two shapes slightly to the right.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Step 1: edits to the existing slide 2 ("Pipeline" diagram)
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$s2.Shapes.Item("Rectangle 21").TextFrame.TextRange.Text = "Handler Pipeline"
$s2.Shapes.Item("Rectangle 21").TextFrame.TextRange.Font.Size = 14

$s2.Shapes.Item("Rectangle 5").TextFrame.TextRange.Text = "Incoming Pipeline"
$s2.Shapes.Item("Rectangle 5").TextFrame.TextRange.Font.Size = 14

$s2.Shapes.Item("Rectangle 7").TextFrame.TextRange.Text = "Outgoing Pipeline"
$s2.Shapes.Item("Rectangle 7").TextFrame.TextRange.Font.Size = 14

$s2.Shapes.Item("Left Brace 9").Left = 6586798 / 12700
$s2.Shapes.Item("TextBox 20").Left = 6888701 / 12700

# ---------------------------------------------------------------------
# Step 2: duplicate slide 2 to create the new slide 3
# ---------------------------------------------------------------------
$dup = $s2.Duplicate()
$s3 = $p.Slides.Item(3)

# -- Handler Pipeline rectangle (top) --
$s3.Shapes.Item("Rectangle 21").Left = 6931293 / 12700
$s3.Shapes.Item("Rectangle 21").Top = 3088595 / 12700

# -- Incoming Pipeline rectangle --
$s3.Shapes.Item("Rectangle 5").Width = 2664092 / 12700

# -- User Code rectangle --
$s3.Shapes.Item("Rectangle 6").Left = 8291297 / 12700
$s3.Shapes.Item("Rectangle 6").Top = 3091544 / 12700

# -- Outgoing Pipeline rectangle --
$s3.Shapes.Item("Rectangle 7").Left = 6931293 / 12700
$s3.Shapes.Item("Rectangle 7").Width = 1385398 / 12700

# -- Left Brace 9: rotate 90 degrees and reposition --
$s3.Shapes.Item("Left Brace 9").Left = 4576204 / 12700
$s3.Shapes.Item("Left Brace 9").Top = 3491425 / 12700
$s3.Shapes.Item("Left Brace 9").Rotation = 90

# -- Group 17 (DB connection star+label) reposition --
$s3.Shapes.Item("Group 17").Left = 5073347 / 12700
$s3.Shapes.Item("Group 17").Top = 4171798 / 12700

# -- TextBox 20: "Handler"/"Scope" -> "Behavior Scope" --
$tb20 = $s3.Shapes.Item("TextBox 20")
$tb20.TextFrame.TextRange.Text = "Behavior Scope"
$tb20.TextFrame.WordWrap = $true
$tb20.Left = 5301145 / 12700
$tb20.Top = 3917122 / 12700
$tb20.Width = 1798060 / 12700
$tb20.Height = 369332 / 12700

# -- Left Brace 23: rotate 90 degrees and reposition --
$s3.Shapes.Item("Left Brace 23").Left = 4575620 / 12700
$s3.Shapes.Item("Left Brace 23").Top = 3201046 / 12700
$s3.Shapes.Item("Left Brace 23").Rotation = 90

# -- TextBox 24: "Suppress"/"Scope" -> "Suppress Scope" --
$tb24 = $s3.Shapes.Item("TextBox 24")
$tb24.TextFrame.TextRange.Text = "Suppress Scope"
$tb24.TextFrame.WordWrap = $true
$tb24.Left = 5301145 / 12700
$tb24.Top = 3651893 / 12700
$tb24.Width = 1687146 / 12700
$tb24.Height = 369332 / 12700

# -- Straight Connector 25: becomes a vertical line --
$conn = $s3.Shapes.Item("Straight Connector 25")
$conn.Left = 6874677 / 12700
$conn.Top = 3755819 / 12700
$conn.Width = 0 / 12700
$conn.Height = 895777 / 12700

# -- New Cylinder shape (DB icon) --
$cyl = $s3.Shapes.AddShape(22, 4583927 / 12700, 4246121 / 12700, 402671 / 12700, 293524 / 12700)
$cyl.Name = "Cylinder 1"
